# Scheduled runner update: refresh cached market-board price data
# (currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfit columns)
# for a handful of leve rows across the per-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 7500
$ws.Range("J74").Value = 7500
$ws.Range("L74").Value = 7500
$ws.Range("N74").Value = -9372

$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()

$ws.Range("H77").Value = 7500
$ws.Range("J77").Value = 7500
$ws.Range("L77").Value = 37500
$ws.Range("N77").Value = -46860

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()

$ws.Range("H88").Value = 3493.258
$ws.Range("I88").Value = 2666.6667
$ws.Range("J88").Value = 3581.8215
$ws.Range("K88").Value = 2666.6667
$ws.Range("L88").Value = 3581.8215
$ws.Range("M88").Value = -2260.6667
$ws.Range("N88").Value = -4393.8215

$ws.Range("H91").Value = 3493.258
$ws.Range("I91").Value = 2666.6667
$ws.Range("J91").Value = 3581.8215
$ws.Range("K91").Value = 2666.6667
$ws.Range("L91").Value = 3581.8215
$ws.Range("M91").Value = -1262.6667
$ws.Range("N91").Value = -6389.8215

$ws.Range("H103").Value = 1125.75
$ws.Range("I103").Value = 1334.6666
$ws.Range("J103").Value = 499
$ws.Range("K103").Value = 4003.9998
$ws.Range("L103").Value = 1497
$ws.Range("M103").Value = -3417.9998
$ws.Range("N103").Value = -2669

$ws.Range("H112").Value = 4499.778
$ws.Range("I112").Value = 758.5
$ws.Range("J112").Value = 4967.4375
$ws.Range("K112").Value = 2275.5
$ws.Range("L112").Value = 14902.3125
$ws.Range("M112").Value = -1167.5
$ws.Range("N112").Value = -17118.3125

$ws.Range("H135").Value = 1106.174
$ws.Range("I135").Value = 508.8889
$ws.Range("K135").Value = 4580.0001
$ws.Range("M135").Value = -2045.0001

$ws.Range("H137").Value = 49252.58
$ws.Range("I137").Value = 63533.45
$ws.Range("J137").Value = 3236.4443
$ws.Range("K137").Value = 190600.35
$ws.Range("L137").Value = 9709.332900000001
$ws.Range("M137").Value = -188050.35
$ws.Range("N137").Value = -14809.3329

$ws.Range("H140").Value = 94115.60000000001
$ws.Range("J140").Value = 94115.60000000001
$ws.Range("L140").Value = 94115.60000000001
$ws.Range("N140").Value = -104475.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 51796.5
$ws.Range("I2").Value = 64142.5
$ws.Range("K2").Value = 64142.5
$ws.Range("M2").Value = -64029.5

$ws.Range("H32").Value = 8834.380999999999
$ws.Range("I32").Value = 5156.656
$ws.Range("K32").Value = 5156.656
$ws.Range("M32").Value = -4869.656

$ws.Range("H45").Value = 71741.60000000001
$ws.Range("I45").Value = 113126.89
$ws.Range("K45").Value = 113126.89
$ws.Range("M45").Value = -112749.89

$ws.Range("H61").Value = 3209.9443
$ws.Range("I61").Value = 3049.4333
$ws.Range("K61").Value = 3049.4333
$ws.Range("M61").Value = -2837.4333

$ws.Range("H74").Value = 57996.47
$ws.Range("I74").Value = 40532
$ws.Range("J74").Value = 102627.89
$ws.Range("K74").Value = 40532
$ws.Range("L74").Value = 102627.89
$ws.Range("M74").Value = -39658
$ws.Range("N74").Value = -104375.89

$ws.Range("H77").Value = 57996.47
$ws.Range("I77").Value = 40532
$ws.Range("J77").Value = 102627.89
$ws.Range("K77").Value = 202660
$ws.Range("L77").Value = 513139.45
$ws.Range("M77").Value = -198292
$ws.Range("N77").Value = -521875.45

$ws.Range("H116").Value = 51796.5
$ws.Range("I116").Value = 64142.5
$ws.Range("K116").Value = 64142.5
$ws.Range("M116").Value = -61848.5

$ws.Range("H122").Value = 9263767
$ws.Range("I122").Value = 14817628
$ws.Range("K122").Value = 44452884
$ws.Range("M122").Value = -44450434

$ws.Range("H132").Value = 3074.111
$ws.Range("I132").Value = 2199.8572
$ws.Range("K132").Value = 6599.571599999999
$ws.Range("M132").Value = -4069.571599999999

$ws.Range("H136").Value = 3209.9443
$ws.Range("I136").Value = 3049.4333
$ws.Range("K136").Value = 9148.2999
$ws.Range("M136").Value = -6598.2999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 51796.5
$ws.Range("I3").Value = 64142.5
$ws.Range("K3").Value = 64142.5
$ws.Range("M3").Value = -64028.5

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

$ws.Range("H80").Value = 408.2195
$ws.Range("I80").Value = 356.16666
$ws.Range("J80").Value = 448.9565
$ws.Range("K80").Value = 356.16666
$ws.Range("L80").Value = 448.9565
$ws.Range("M80").Value = 641.83334
$ws.Range("N80").Value = -2444.9565

$ws.Range("H83").Value = 408.2195
$ws.Range("I83").Value = 356.16666
$ws.Range("J83").Value = 448.9565
$ws.Range("K83").Value = 1780.8333
$ws.Range("L83").Value = 2244.7825
$ws.Range("M83").Value = 3211.1667
$ws.Range("N83").Value = -12228.7825

$ws.Range("H86").Value = 5402.857
$ws.Range("I86").Value = 7035.737
$ws.Range("J86").Value = 1955.6666
$ws.Range("K86").Value = 7035.737
$ws.Range("L86").Value = 1955.6666
$ws.Range("M86").Value = -5912.737
$ws.Range("N86").Value = -4201.6666

$ws.Range("H89").Value = 5402.857
$ws.Range("I89").Value = 7035.737
$ws.Range("J89").Value = 1955.6666
$ws.Range("K89").Value = 35178.685
$ws.Range("L89").Value = 9778.333000000001
$ws.Range("M89").Value = -29562.685
$ws.Range("N89").Value = -21010.333

$ws.Range("H99").Value = 3003.3157
$ws.Range("I99").Value = 2391.4614
$ws.Range("K99").Value = 2391.4614
$ws.Range("M99").Value = -893.4614000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26202.309
$ws.Range("I31").Value = 1621.9546
$ws.Range("J31").Value = 58012.176
$ws.Range("K31").Value = 1621.9546
$ws.Range("L31").Value = 58012.176
$ws.Range("M31").Value = -1326.9546
$ws.Range("N31").Value = -58602.176

$ws.Range("H34").Value = 26202.309
$ws.Range("I34").Value = 1621.9546
$ws.Range("J34").Value = 58012.176
$ws.Range("K34").Value = 1621.9546
$ws.Range("L34").Value = 58012.176
$ws.Range("M34").Value = -1419.9546
$ws.Range("N34").Value = -58416.176

$ws.Range("H105").Value = 1187.3334
$ws.Range("I105").Value = 1187.3334
$ws.Range("K105").Value = 1187.3334
$ws.Range("M105").Value = 559.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 168
$ws.Range("I109").Value = 168
$ws.Range("K109").Value = 504
$ws.Range("M109").Value = 536

$ws.Range("H131").Value = 12629387
$ws.Range("J131").Value = 11908279
$ws.Range("L131").Value = 35724837
$ws.Range("N131").Value = -35734917

$ws.Range("H137").Value = 2071.5
$ws.Range("I137").Value = 1777.6666
$ws.Range("J137").Value = 2600.4
$ws.Range("K137").Value = 5332.9998
$ws.Range("L137").Value = 7801.200000000001
$ws.Range("M137").Value = -232.9997999999996
$ws.Range("N137").Value = -18001.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 71444510
$ws.Range("I80").Value = 83351100
$ws.Range("K80").Value = 83351100
$ws.Range("M80").Value = -83350102

$ws.Range("H83").Value = 71444510
$ws.Range("I83").Value = 83351100
$ws.Range("K83").Value = 416755500
$ws.Range("M83").Value = -416750508

$ws.Range("H141").Value = 42567.2
$ws.Range("I141").Value = 25000
$ws.Range("J141").Value = 46959
$ws.Range("K141").Value = 25000
$ws.Range("L141").Value = 46959
$ws.Range("N141").Value = -57319
$ws.Range("M141").Value = -19820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 34818.168
$ws.Range("I4").Value = 32449.5
$ws.Range("J4").Value = 36002.5
$ws.Range("K4").Value = 32449.5
$ws.Range("L4").Value = 36002.5
$ws.Range("M4").Value = -32336.5
$ws.Range("N4").Value = -36228.5

$ws.Range("H28").Value = 34818.168
$ws.Range("I28").Value = 32449.5
$ws.Range("J28").Value = 36002.5
$ws.Range("K28").Value = 32449.5
$ws.Range("L28").Value = 36002.5
$ws.Range("M28").Value = -32217.5
$ws.Range("N28").Value = -36466.5

$ws.Range("H37").Value = 34818.168
$ws.Range("I37").Value = 32449.5
$ws.Range("J37").Value = 36002.5
$ws.Range("K37").Value = 32449.5
$ws.Range("L37").Value = 36002.5
$ws.Range("M37").Value = -32342.5
$ws.Range("N37").Value = -36216.5

$ws.Range("H127").Value = 59237.5
$ws.Range("J127").Value = 59237.5
$ws.Range("L127").Value = 59237.5
$ws.Range("N127").Value = -69157.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 64999
$ws.Range("J141").Value = 64999
$ws.Range("L141").Value = 64999
$ws.Range("N141").Value = -75359
